# Actualización automática 2025-07-09 12:25:08
# Registers a new "julio" (July) sale for two clients of asesor
# HIDALGO HIDALGO PEDRO GUSTAVO and propagates the totals across the
# three sheets of the workbook.

$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo    = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual      = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento      = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": per-client amount sold per product group
# ---------------------------------------------------------------------
# CHASIQUIZA CAMPAÑA JOSE LUIS (row 6) -> 240X80 PORCELANATO (col D)
$wsVentasPorGrupo.Range("D6").Value = 2472.77
# MUÑOZ LOZA ROMMEL SEBASTIAN (row 13) -> PORCELANATO (col M)
$wsVentasPorGrupo.Range("M13").Value = 1451.52

# Row 22 keeps a "N de 20" count of clients with sales per group; the
# new sale pushes the count for 240X80 PORCELANATO and the overall
# PORCELANATO total up by one.
$wsVentasPorGrupo.Range("D22").Value = "2 de 20"
$wsVentasPorGrupo.Range("M22").Value = "6 de 20"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": per-client sales by month
# ---------------------------------------------------------------------
# julio (July) is column F
$wsVentaMensual.Range("F6").Value = 2472.77
$wsVentaMensual.Range("F13").Value = 1451.52
$wsVentaMensual.Range("F22").Value = 23227.39

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": budget vs sales per product group
# ---------------------------------------------------------------------
# Row 3: 240X80 PORCELANATO
$wsCumplimiento.Range("D3").Value = 2915.04
$wsCumplimiento.Range("E3").Value = 1253.03156573679
$wsCumplimiento.Range("F3").Value = 0.6993737880997032

# Row 16: PORCELANATO
$wsCumplimiento.Range("D16").Value = 17398.67
$wsCumplimiento.Range("E16").Value = 26867.57
$wsCumplimiento.Range("F16").Value = 0.3930460323713963

# Row 19: TOTAL
$wsCumplimiento.Range("D19").Value = 23227.39
$wsCumplimiento.Range("E19").Value = 42150.60762291768
$wsCumplimiento.Range("F19").Value = 0.3552783940243811
